# Add the missing constraints: update the computed "Value:" column (C) for
# several variables on the "Var Type 1", "Var Type 2" and "Var Type AFT"
# sheets, reflecting the corrected constraints in the AffineFlowThinning
# script.

$wb = $excel.ActiveWorkbook

$wsType1 = $wb.Worksheets.Item("Var Type 1")
$wsType1.Range("C2").Value = 10
$wsType1.Range("C3").Value = 5
$wsType1.Range("C4").Value = 8.75

$wsType2 = $wb.Worksheets.Item("Var Type 2")
$wsType2.Range("C2").Value = 5
$wsType2.Range("C4").Value = 7.75
$wsType2.Range("C8").Value = 3

$wsTypeAFT = $wb.Worksheets.Item("Var Type AFT")
$wsTypeAFT.Range("C2").Value = 10
$wsTypeAFT.Range("C5").Value = 8.75
$wsTypeAFT.Range("C11").Value = 5
